$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.12
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4.2
$ws.Range("L2").Value = 1.32
$ws.Range("Q2").Value = 1.7
$ws.Range("R2").Value = 1.51
$ws.Range("Z2").Value = 16.5
$ws.Range("AA2").Value = 28
$ws.Range("AD2").Value = 11.5
$ws.Range("AE2").Value = 22
$ws.Range("AF2").Value = 28
$ws.Range("AH2").Value = 18.5
$ws.Range("AL2").Value = 42
$ws.Range("AO2").Value = 13.5
$ws.Range("P3").Value = 2.22
$ws.Range("Q3").Value = 1.71
$ws.Range("R3").Value = 1.47
$ws.Range("S3").Value = 2.84
$ws.Range("U3").Value = 2.36
$ws.Range("AN3").Value = 13.5
$ws.Range("F4").Value = 1.66
$ws.Range("G4").Value = 1.78
$ws.Range("H4").Value = 5.1
$ws.Range("I4").Value = 6.8
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4.6
$ws.Range("P4").Value = 1.98
$ws.Range("Q4").Value = 1.86
$ws.Range("F5").Value = 3.15
$ws.Range("I5").Value = 3.25
$ws.Range("G6").Value = 2.4
$ws.Range("J6").Value = 2.98
$ws.Range("F7").Value = 1.9
$ws.Range("F8").Value = 1.63
$ws.Range("G8").Value = 1.77
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 4.3
$ws.Range("K8").Value = 5.2
$ws.Range("P8").Value = 1.25
$ws.Range("Q8").Value = 1.01
$ws.Range("P9").Value = 1.25
$ws.Range("H10").Value = 3
$ws.Range("G11").Value = 1.77
$ws.Range("I11").Value = 8.4
$ws.Range("J11").Value = 3.65
$ws.Range("P11").Value = 1.83
$ws.Range("F12").Value = 2.18
$ws.Range("P12").Value = 1.98
$ws.Range("U12").Value = 2.18
$ws.Range("X12").Value = 14.5
$ws.Range("P13").Value = 2.12
$ws.Range("X13").Value = 17.5
$ws.Range("H14").Value = 3.65
$ws.Range("N14").Value = 2.92
$ws.Range("AH14").Value = 22
$ws.Range("AK14").Value = 34
$ws.Range("T16").Value = 2.04
$ws.Range("G17").Value = 3.8
$ws.Range("P17").Value = 1.45
$ws.Range("Q17").Value = 2.88
